$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Igf2r"
$ws.Cells.Item(2,3).Value = "ECs"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 13.24090133333333
$ws.Cells.Item(2,8).Value = 39.722704
$ws.Cells.Item(2,9).Value = 0.1214410874295642
$ws.Cells.Item(2,10).Value = 0.1214410874295642
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 36.51516
$ws.Cells.Item(2,14).Value = 109.54548
$ws.Cells.Item(2,15).Value = 0.3318909895952502
$ws.Cells.Item(2,16).Value = 0.3318909895952502
$ws.Cells.Item(2,17).Value = 483.49363073088
$ws.Cells.Item(2,18).Value = 4351.44267657792
$ws.Cells.Item(2,19).Value = 0.04030520268452136
$ws.Cells.Item(2,20).Value = 0.04030520268452137

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Igf2r"
$ws.Cells.Item(3,3).Value = "ECs"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 13.24090133333333
$ws.Cells.Item(3,8).Value = 39.722704
$ws.Cells.Item(3,9).Value = 0.1214410874295642
$ws.Cells.Item(3,10).Value = 0.1214410874295642
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 42.26455300000001
$ws.Cells.Item(3,14).Value = 126.793659
$ws.Cells.Item(3,15).Value = 0.3841479626536184
$ws.Cells.Item(3,16).Value = 0.3841479626536184
$ws.Cells.Item(3,17).Value = 559.6207761704375
$ws.Cells.Item(3,18).Value = 5036.586985533937
$ws.Cells.Item(3,19).Value = 0.04665134631850704
$ws.Cells.Item(3,20).Value = 0.04665134631850705

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Igf2r"
$ws.Cells.Item(4,3).Value = "ECs"
$ws.Cells.Item(4,4).Value = "Igf2"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 13.24090133333333
$ws.Cells.Item(4,8).Value = 39.722704
$ws.Cells.Item(4,9).Value = 0.1214410874295642
$ws.Cells.Item(4,10).Value = 0.1214410874295642
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 31.24183366666666
$ws.Cells.Item(4,14).Value = 93.725501
$ws.Cells.Item(4,15).Value = 0.2839610477511313
$ws.Cells.Item(4,16).Value = 0.2839610477511314
$ws.Cells.Item(4,17).Value = 413.6700370527449
$ws.Cells.Item(4,18).Value = 3723.030333474704
$ws.Cells.Item(4,19).Value = 0.0344845384265358
$ws.Cells.Item(4,20).Value = 0.03448453842653581

# Row 5
$ws.Cells.Item(5,1).Value = "sCs"
$ws.Cells.Item(5,2).Value = "Igf2r"
$ws.Cells.Item(5,3).Value = "ECs"
$ws.Cells.Item(5,4).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 5.270503666666666
$ws.Cells.Item(5,8).Value = 15.811511
$ws.Cells.Item(5,9).Value = 0.04833928449947708
$ws.Cells.Item(5,10).Value = 0.04833928449947708
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 36.51516
$ws.Cells.Item(5,14).Value = 109.54548
$ws.Cells.Item(5,15).Value = 0.3318909895952502
$ws.Cells.Item(5,16).Value = 0.3318909895952502
$ws.Cells.Item(5,17).Value = 192.45328466892
$ws.Cells.Item(5,18).Value = 1732.07956202028
$ws.Cells.Item(5,19).Value = 0.01604337296885779
$ws.Cells.Item(5,20).Value = 0.01604337296885779

# Row 6
$ws.Cells.Item(6,1).Value = "sCs"
$ws.Cells.Item(6,2).Value = "Igf2r"
$ws.Cells.Item(6,3).Value = "ECs"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 5.270503666666666
$ws.Cells.Item(6,8).Value = 15.811511
$ws.Cells.Item(6,9).Value = 0.04833928449947708
$ws.Cells.Item(6,10).Value = 0.04833928449947708
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 42.26455300000001
$ws.Cells.Item(6,14).Value = 126.793659
$ws.Cells.Item(6,15).Value = 0.3841479626536184
$ws.Cells.Item(6,16).Value = 0.3841479626536184
$ws.Cells.Item(6,17).Value = 222.7554815565277
$ws.Cells.Item(6,18).Value = 2004.799334008749
$ws.Cells.Item(6,19).Value = 0.01856943765660776
$ws.Cells.Item(6,20).Value = 0.01856943765660776

# Row 7
$ws.Cells.Item(7,1).Value = "sCs"
$ws.Cells.Item(7,2).Value = "Igf2r"
$ws.Cells.Item(7,3).Value = "ECs"
$ws.Cells.Item(7,4).Value = "Igf2"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 5.270503666666666
$ws.Cells.Item(7,8).Value = 15.811511
$ws.Cells.Item(7,9).Value = 0.04833928449947708
$ws.Cells.Item(7,10).Value = 0.04833928449947708
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 31.24183366666666
$ws.Cells.Item(7,14).Value = 93.725501
$ws.Cells.Item(7,15).Value = 0.2839610477511313
$ws.Cells.Item(7,16).Value = 0.2839610477511314
$ws.Cells.Item(7,17).Value = 164.6601988935568
$ws.Cells.Item(7,18).Value = 1481.941790042011
$ws.Cells.Item(7,19).Value = 0.01372647387401153
$ws.Cells.Item(7,20).Value = 0.01372647387401154

# Row 8
$ws.Cells.Item(8,1).Value = "Igf2"
$ws.Cells.Item(8,2).Value = "Igf2r"
$ws.Cells.Item(8,3).Value = "ECs"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 90.52007366666668
$ws.Cells.Item(8,8).Value = 271.560221
$ws.Cells.Item(8,9).Value = 0.8302196280709587
$ws.Cells.Item(8,10).Value = 0.8302196280709586
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 36.51516
$ws.Cells.Item(8,14).Value = 109.54548
$ws.Cells.Item(8,15).Value = 0.3318909895952502
$ws.Cells.Item(8,16).Value = 0.3318909895952502
$ws.Cells.Item(8,17).Value = 3305.35497315012
$ws.Cells.Item(8,18).Value = 29748.19475835108
$ws.Cells.Item(8,19).Value = 0.275542413941871
$ws.Cells.Item(8,20).Value = 0.275542413941871

# Row 9
$ws.Cells.Item(9,1).Value = "Igf2"
$ws.Cells.Item(9,2).Value = "Igf2r"
$ws.Cells.Item(9,3).Value = "ECs"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 90.52007366666668
$ws.Cells.Item(9,8).Value = 271.560221
$ws.Cells.Item(9,9).Value = 0.8302196280709587
$ws.Cells.Item(9,10).Value = 0.8302196280709586
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 42.26455300000001
$ws.Cells.Item(9,14).Value = 126.793659
$ws.Cells.Item(9,15).Value = 0.3841479626536184
$ws.Cells.Item(9,16).Value = 0.3841479626536184
$ws.Cells.Item(9,17).Value = 3825.790451048739
$ws.Cells.Item(9,18).Value = 34432.11405943865
$ws.Cells.Item(9,19).Value = 0.3189271786785036
$ws.Cells.Item(9,20).Value = 0.3189271786785036

# Row 10
$ws.Cells.Item(10,1).Value = "Igf2"
$ws.Cells.Item(10,2).Value = "Igf2r"
$ws.Cells.Item(10,3).Value = "ECs"
$ws.Cells.Item(10,4).Value = "Igf2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 90.52007366666668
$ws.Cells.Item(10,8).Value = 271.560221
$ws.Cells.Item(10,9).Value = 0.8302196280709587
$ws.Cells.Item(10,10).Value = 0.8302196280709586
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 31.24183366666666
$ws.Cells.Item(10,14).Value = 93.725501
$ws.Cells.Item(10,15).Value = 0.2839610477511313
$ws.Cells.Item(10,16).Value = 0.2839610477511314
$ws.Cells.Item(10,17).Value = 2828.013084988414
$ws.Cells.Item(10,18).Value = 25452.11776489572
$ws.Cells.Item(10,19).Value = 0.235750035450584
$ws.Cells.Item(10,20).Value = 0.235750035450584
